$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 14:05"

# Country name re-sort (column A) caused by updated data ordering
$ws.Range("A37").Value = "Kuwait"
$ws.Range("A38").Value = "Sudafrica"
$ws.Range("A39").Value = "Colombia"
$ws.Range("A40").Value = "Rumania"
$ws.Range("A105").Value = "Sri Lanka"
$ws.Range("A106").Value = "Kenia"
$ws.Range("A108").Value = "Libano"
$ws.Range("A109").Value = "Albania"
$ws.Range("A124").Value = "Haiti"
$ws.Range("A125").Value = "San Marino"
$ws.Range("A126").Value = "Malta"

# Updated numeric statistics
$ws.Range("B6").Value = 294152
$ws.Range("C6").Value = 795
$ws.Range("E6").Value = 158431
$ws.Range("G6").Value = 144
$ws.Range("H6").Value = 19038
$ws.Range("B37").Value = 18609
$ws.Range("C37").Value = 1041
$ws.Range("D37").Value = 5205
$ws.Range("E37").Value = 13275
$ws.Range("G37").Value = 5
$ws.Range("H37").Value = 129
$ws.Range("B38").Value = 18003
$ws.Range("D38").Value = 8950
$ws.Range("E38").Value = 8714
$ws.Range("H38").Value = 339
$ws.Range("B39").Value = 17687
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 4256
$ws.Range("E39").Value = 12801
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 630
$ws.Range("B40").Value = 17585
$ws.Range("C40").Value = 198
$ws.Range("D40").Value = 10581
$ws.Range("E40").Value = 5853
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 1151
$ws.Range("D76").Value = 2407
$ws.Range("E76").Value = 530
$ws.Range("B105").Value = 1030
$ws.Range("C105").Value = 2
$ws.Range("D105").Value = 604
$ws.Range("E105").Value = 417
$ws.Range("H105").Value = 9
$ws.Range("B106").Value = 1029
$ws.Range("D106").Value = 366
$ws.Range("E106").Value = 613
$ws.Range("H106").Value = 50
$ws.Range("B108").Value = 1024
$ws.Range("C108").Value = 63
$ws.Range("D108").Value = 251
$ws.Range("E108").Value = 747
$ws.Range("H108").Value = 26
$ws.Range("B109").Value = 969
$ws.Range("C109").Value = 5
$ws.Range("D109").Value = 771
$ws.Range("E109").Value = 167
$ws.Range("H109").Value = 31
$ws.Range("B124").Value = 663
$ws.Range("C124").Value = 67
$ws.Range("D124").Value = 21
$ws.Range("E124").Value = 620
$ws.Range("H124").Value = 22
$ws.Range("B125").Value = 658
$ws.Range("C125").Value = 2
$ws.Range("D125").Value = 235
$ws.Range("E125").Value = 382
$ws.Range("H125").Value = 41
$ws.Range("B126").Value = 599
$ws.Range("C126").Value = 15
$ws.Range("D126").Value = 468
$ws.Range("E126").Value = 125
$ws.Range("H126").Value = 6
